$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ансимов Артём's phone number was stored as a raw number; store as text
# (same digits, now a shared string) - set this first so shared-string
# ordering matches the source edit.
$ws.Range("B4").Value = "89834626819"

# Гамаюнова Аделина's phone number now has a leading "+7" instead of "8"
$ws.Range("B2").Value = "+79842740104"

# Column A was widened to fit the longer name currently shown in the UI
$ws.Columns.Item(1).ColumnWidth = 26

# Active cell / selection moved to B3
$ws.Range("B3").Select()
